# Adds <w:proofErr> spell-check bracket markers around a handful of runs
# ("M", "Healthineers", "Springboot", "Github" x3, "M") by splitting the
# runs that contained them and wrapping the relevant run with
# w:proofErr[type=spellStart] / w:proofErr[type=spellEnd], exactly as the
# OOXML diff describes. Word itself stamps these in as a side effect of
# its background spell-checker; since this runtime has no live spell
# checker we reconstruct each affected paragraph's WordprocessingML by
# hand and splice it back in with Range.InsertXML (which, for this
# engine, replaces the *entire* paragraph that the target Range sits
# in - so every replacement below carries the paragraph's full original
# <w:pPr> and the full run list, with only the intended run(s) split).

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml {
    param(
        [string]$FindText,
        [string]$InnerXml
    )

    $rng = $d.Content
    $ok = $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $FindText"
    }
    $para = $rng.Paragraphs(1)
    $prange = $para.Range

    $package = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document $wns>
<w:body>
$InnerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

    $prange.InsertXML($package)
}

# 1) Title: "Srikara" + " M M"  ->  "Srikara" + " M " + proofErr(spellStart) "M" proofErr(spellEnd)
Set-ParagraphXml -FindText "Srikara" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="6128"/>
      <w:tab w:val="left" w:pos="6296"/>
      <w:tab w:val="left" w:pos="6396"/>
      <w:tab w:val="left" w:pos="6513"/>
    </w:tabs>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>Srikara</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> M </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>M</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

# 2) "Medical Healthineers" -> "Medical " + proofErr(spellStart) "Healthineers" proofErr(spellEnd)
Set-ParagraphXml -FindText "Medical Healthineers" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="30"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="1440"/>
      <w:tab w:val="left" w:pos="2629"/>
    </w:tabs>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Medical </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Healthineers</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

# 3) "Springboot" (table cell) -> proofErr(spellStart) "Springboot" proofErr(spellEnd) (no text change)
Set-ParagraphXml -FindText "Springboot" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="34"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="630"/>
      <w:tab w:val="left" w:pos="810"/>
      <w:tab w:val="left" w:pos="900"/>
      <w:tab w:val="left" w:pos="3240"/>
      <w:tab w:val="left" w:pos="3420"/>
      <w:tab w:val="left" w:pos="3510"/>
    </w:tabs>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Springboot</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

# 4) "Github & Github Desktop" ->
#    proofErr(spellStart) "Github" proofErr(spellEnd) + " &" + " " +
#    proofErr(spellStart) "Github" proofErr(spellEnd) + " Desktop"
Set-ParagraphXml -FindText "Github & Github Desktop" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="32"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="630"/>
      <w:tab w:val="left" w:pos="720"/>
      <w:tab w:val="left" w:pos="810"/>
      <w:tab w:val="left" w:pos="900"/>
      <w:tab w:val="left" w:pos="3240"/>
      <w:tab w:val="left" w:pos="3420"/>
      <w:tab w:val="left" w:pos="3510"/>
    </w:tabs>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Github</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> &amp;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Github</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Desktop</w:t>
  </w:r>
</w:p>
'@

# 5) "Master Git & Github" -> "Master Git & " + proofErr(spellStart) "Github" proofErr(spellEnd)
Set-ParagraphXml -FindText "Master Git & Github" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="32"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="630"/>
      <w:tab w:val="left" w:pos="720"/>
      <w:tab w:val="left" w:pos="810"/>
      <w:tab w:val="left" w:pos="900"/>
      <w:tab w:val="left" w:pos="3240"/>
      <w:tab w:val="left" w:pos="3420"/>
      <w:tab w:val="left" w:pos="3510"/>
    </w:tabs>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Master Git &amp; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Github</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

# 6) Footer: "...M M SRIKARA" -> "...M " + proofErr(spellStart) "M" proofErr(spellEnd) + " SRIKARA"
Set-ParagraphXml -FindText "M M SRIKARA" -InnerXml @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="426"/>
      <w:tab w:val="left" w:pos="810"/>
      <w:tab w:val="left" w:pos="900"/>
      <w:tab w:val="left" w:pos="3240"/>
      <w:tab w:val="left" w:pos="3420"/>
      <w:tab w:val="left" w:pos="3510"/>
    </w:tabs>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:ind w:left="567" w:hanging="578"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>Place : Bengaluru</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">     </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">        </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve">M </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>M</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve"> SRIKARA</w:t>
  </w:r>
</w:p>
'@

Write-Output "done"
